$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 85.8724807945396

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
$ws.Range("N6").Value = $newValue
